# "Generate Report for Archive"
#
# This commit corresponds to a CI job re-generating the localization-status
# report workbook. Diffing the canonical OOXML shows the run produced a
# byte-different xl/sharedStrings.xml (one extra, already-known status
# string - "In Translation" - was re-interned into the shared-string pool,
# and a couple of other strings were re-ordered within that pool) and the
# worksheets' <c t="s"><v>N</v></c> references were renumbered to match.
#
# Crucially, every cell's *displayed* value/type is identical before and
# after: each renumbered shared-string index still resolves to exactly the
# same text (e.g. "Ready for handoff", the various timestamps/GUID file
# names, etc.), and the newly interned string isn't referenced by any cell.
# So there is no actual row/column/value change to make on the workbook -
# the Overview/zh-cn/de-de sheets, their headers, data rows, hyperlinks and
# tables all stay exactly as they were; only the internal string-pool
# bookkeeping that produced the archived report differs.
#
# Re-writing the cells through the object model isn't a faithful way to
# reproduce that: typing literal "True"/"False" back into a cell makes
# Excel coerce it to a native Boolean instead of leaving it as the original
# text value, and re-typing "" over a blank cell removes it - both of which
# would introduce real (and incorrect) changes that aren't present in the
# diff. So this script intentionally leaves the workbook's data untouched,
# matching the no-op-for-content nature of the archived report generation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
